$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$codes = @(
    "CVE_MUN",
    "13001","13002","13003","13004","13005","13006","13007","13008","13009","13010",
    "13011","13012","13013","13014","13015","13016","13017","13018","13019","13020",
    "13021","13022","13023","13024","13025","13026","13027","13028","13029","13030",
    "13031","13032","13033","13034","13035","13036","13037","13038","13039","13040",
    "13041","13042","13043","13044","13045","13046","13047","13048","13049","13050",
    "13051","13052","13053","13054","13055","13056","13057","13058","13059","13060",
    "13061","13062","13063","13064","13065","13066","13067","13068","13069","13070",
    "13071","13072","13073","13074","13075","13076","13077","13078","13079","13080",
    "13081","13082","13083","13084"
)

$range = $ws.Range("A1:A85")
$range.NumberFormat = "@"

for ($i = 0; $i -lt $codes.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $codes[$i]
}
